$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-08-28 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-29 Tuesday", 2) | Out-Null

# Update the multiplication problems in the table, addressed by (row, column)
# to avoid ambiguity since some new values equal other cells old values.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "96×74="  # was "39×44="
$t.Cell(1, 2).Range.Text = "88×27="  # was "28×67="
$t.Cell(1, 3).Range.Text = "46×18="  # was "95×20="
$t.Cell(1, 4).Range.Text = "91×77="  # was "63×18="
$t.Cell(1, 5).Range.Text = "59×84="  # was "97×18="
$t.Cell(5, 1).Range.Text = "95×71="  # was "57×14="
$t.Cell(5, 2).Range.Text = "60×59="  # was "76×63="
$t.Cell(5, 3).Range.Text = "76×58="  # was "29×21="
$t.Cell(5, 4).Range.Text = "12×43="  # was "98×13="
$t.Cell(5, 5).Range.Text = "19×97="  # was "62×54="
$t.Cell(10, 1).Range.Text = "21×19="  # was "72×88="
$t.Cell(10, 2).Range.Text = "95×51="  # was "19×13="
$t.Cell(10, 3).Range.Text = "84×14="  # was "81×41="
$t.Cell(10, 4).Range.Text = "27×58="  # was "67×93="
$t.Cell(10, 5).Range.Text = "76×82="  # was "19×95="
$t.Cell(15, 1).Range.Text = "92×96="  # was "80×58="
$t.Cell(15, 2).Range.Text = "18×13="  # was "89×66="
$t.Cell(15, 3).Range.Text = "83×85="  # was "88×27="
$t.Cell(15, 4).Range.Text = "34×50="  # was "86×29="
$t.Cell(15, 5).Range.Text = "16×66="  # was "63×60="
$t.Cell(20, 1).Range.Text = "75×29="  # was "41×62="
$t.Cell(20, 2).Range.Text = "71×97="  # was "83×22="
$t.Cell(20, 3).Range.Text = "14×87="  # was "60×40="
$t.Cell(20, 4).Range.Text = "50×50="  # was "65×20="
$t.Cell(20, 5).Range.Text = "42×66="  # was "96×44="
